# Append two new GSW @ HOU playoff game rows (row 4 and row 5) to Sheet1,
# matching the existing table layout (header in row 1, data starting row 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: GSW @ HOU (away), 2025-04-23, loss
$row4 = @(2, "GSW", "HOU", "away", "2025-04-23", "240:00", 33, 80, 0.413, 17, 43, 0.395, 11, 18, 0.611, 7, 26, 33, 22, 7, 6, 15, 18, 94, -15, 18, 28, 25, 23, "L")

# Row 5: HOU vs GSW (home), 2025-04-23, win
$row5 = @(3, "HOU", "GSW", "home", "2025-04-23", "240:00", 39, 86, 0.453, 15, 40, 0.375, 16, 20, 0.8, 11, 36, 47, 21, 9, 4, 10, 17, 109, 15, 28, 32, 27, 22, "W")

# Column E ("DATE") holds date-like text (e.g. "2025-04-23") that Excel's
# normal input parsing would silently convert into a date serial number
# (and stamp a date NumberFormat on the cell). Skip it in the generic loop
# below and write it separately via a text-literal formula instead, which
# is then collapsed back down to a plain static value - ending up as
# ordinary text with no formula and no stray number-format.
for ($i = 0; $i -lt $row4.Length; $i++) {
    $col = $i + 1
    if ($col -eq 5) { continue }
    $ws.Cells.Item(4, $col).Value = $row4[$i]
    $ws.Cells.Item(5, $col).Value = $row5[$i]
}

$ws.Range("E4").Formula = "=""2025-04-23"""
$ws.Range("E5").Formula = "=""2025-04-23"""
$ws.Range("E4:E5").Copy() | Out-Null
$ws.Range("E4:E5").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false

# Column A carries the bordered/bold "index" style used by the existing rows
# (row 3's A3 cell already has it) - copy just the formatting, not the value.
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null
$ws.Range("A5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

Write-Host "Appended rows 4 and 5"
